$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q4" quarter sheet by duplicating the existing
#    "2022-Q3" sheet (so it keeps the same column layout / styling), placing
#    it right before "2022-Q3" in the tab order, then rename it and fill in
#    the new quarter's figures.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)

$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

$q4.Range("D2").Value = "'0.24"
$q4.Range("E2").Value = "'68.42"
$q4.Range("F2").Value = "'4.19"
$q4.Range("G2").Value = "'0.0101"

# ---------------------------------------------------------------------------
# 2. Update the "总计" (overview) sheet: shift every quarter label/row down
#    by one and insert the new 2022-Q4 entry at the top, extending the table
#    with an extra trailing row so the oldest quarter (2021-Q2) keeps its
#    place.
# ---------------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")

# Grow the table by copying the formatting of the last existing row into the
# new row 8 first (so styles match), then overwrite every row's values.
$zj.Range("A7").Copy()
$zj.Range("A8").PasteSpecial(-4122)

$zj.Range("B2").Value = "2022-Q4"
$zj.Range("C2").Value = 1
$zj.Range("D2").Value = 0.01

$zj.Range("B3").Value = "2022-Q3"
$zj.Range("C3").Value = 1
$zj.Range("D3").Value = 0.01

$zj.Range("B4").Value = "2022-Q2"
$zj.Range("C4").Value = 1
$zj.Range("D4").Value = 0.01

$zj.Range("B5").Value = "2022-Q1"
$zj.Range("C5").Value = 1
$zj.Range("D5").Value = 0.03

$zj.Range("B6").Value = "2021-Q4"
$zj.Range("C6").Value = 1
$zj.Range("D6").Value = 0.03

$zj.Range("B7").Value = "2021-Q3"
$zj.Range("C7").Value = 1
$zj.Range("D7").Value = 0.03

$zj.Range("A8").Value = 6
$zj.Range("B8").Value = "2021-Q2"
$zj.Range("C8").Value = 1
$zj.Range("D8").Value = 0.03
